$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part-list text clean-ups / MPN corrections ---

# C12, C19, C22, C25, C27  ->  C12,C19,C22,C25,C27  (remove spaces after commas)
$ws.Range("F4").Value = "C12,C19,C22,C25,C27"

# C20's MPN was wrong (GRM022R60J103KE19L) -> corrected part number
$ws.Range("I9").Value = "01016D103KAT2A"

# C1, C2, C3, C15, C18  ->  C1,C2,C3,C15,C18  (remove spaces after commas)
$ws.Range("F12").Value = "C1,C2,C3,C15,C18"

# 1uF cap MPN corrected
$ws.Range("I12").Value = "JMK063ABJ105KP-F"

# RGB LED MPN corrected / expanded
$ws.Range("I29").Value = "LRTBR48G-P9Q7-1+R7S5-26+NP-68"

# --- Restore the view so the sheet is scrolled/selected the same way the
#     author left it (fixing RX/TX lines meant looking at the connector
#     pins around I30) ---
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 18
$ws.Range("I30").Select()
